$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the embedded SQL queries (StatQuery in C2, TabQuery in B2:B7).
#    The join conditions were switched from the generic "id" column names
#    to the fully-qualified "study_id" / "participant_id" column names.
# ---------------------------------------------------------------------------
function Update-Query([string]$cellRef) {
    $cell = $ws.Range($cellRef)
    $text = $cell.Value2
    if ($text -eq $null) { return }

    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $cell.Value = $text
}

Update-Query "C2"
Update-Query "B2"
Update-Query "B3"
Update-Query "B4"
Update-Query "B5"
Update-Query "B6"
Update-Query "B7"

# ---------------------------------------------------------------------------
# 2) Scroll the sheet view back to the top-left (A1) instead of A6.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 3) Widen column C and drop the "best fit" auto-sizing in favor of an
#    explicit custom width.
# ---------------------------------------------------------------------------
$colC = $ws.Columns.Item(3)
$colC.ColumnWidth = 72.8

Write-Output "Edit applied successfully"
